# ExpenseTemplate.xlsx maintenance update:
#  - MasrafMerkezi (cost center) list refreshed from source file:
#      * "5191 67 Burda Animal" cost-center row removed
#      * "5110 Migros  Animal " row duplicated (new row inserted under it)
#  - MasrafKalemi (expense item) "Araç Kiralama" renamed to
#    "Araç Kiralama / Nakliye"

$wb = $excel.ActiveWorkbook

# --- MasrafMerkezi sheet: remove the "5191 67 Burda Animal" cost center
#     row (ID 160) and duplicate the "5110 Migros Animal" row (ID 105) ---
$wsMerkezi = $wb.Worksheets.Item("MasrafMerkezi")

# Row 21 holds ID 160 / "5191 67 Burda Animal" / "Animal" -> delete it.
$wsMerkezi.Rows.Item(21).Delete()

# Insert a new blank row right after row 15 (the "5110 Migros Animal" row)
# and copy that row's values into it, producing a duplicate entry.
$wsMerkezi.Rows.Item(16).Insert()
$wsMerkezi.Range("A16").Value2 = $wsMerkezi.Range("A15").Value2
$wsMerkezi.Range("B16").Value2 = $wsMerkezi.Range("B15").Value2
$wsMerkezi.Range("C16").Value2 = $wsMerkezi.Range("C15").Value2

# --- MasrafKalemi sheet: rename "Araç Kiralama" expense item ---
$wsKalemi = $wb.Worksheets.Item("MasrafKalemi")
$wsKalemi.Range("B22").Value2 = "Araç Kiralama / Nakliye"
